$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 45) captured on 2025-07-27.
# Force the date column to be stored as plain text (matching the existing
# "MM/DD/YYYY" text entries already in column A) instead of letting Excel
# auto-convert the date-like string into a date serial number.
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "07/27/2025"
$ws.Range("A45").Style = "Normal"

$ws.Range("B45").Value = 93.56999999999971
$ws.Range("C45").Value = 0.1068718606390941
$ws.Range("D45").Value = 10
